$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-21 Sunday" "2024-07-22 Monday"

Replace-Text "559÷6=93, 1" "551÷8=68, 7"
Replace-Text "596÷7=85, 1" "155÷5=31, 0"
Replace-Text "360÷2=180, 0" "151÷4=37, 3"
Replace-Text "159÷7=22, 5" "654÷5=130, 4"
Replace-Text "249÷2=124, 1" "939÷2=469, 1"

Replace-Text "106÷9=11, 7" "192÷9=21, 3"
Replace-Text "233÷4=58, 1" "894÷7=127, 5"
Replace-Text "404÷3=134, 2" "524÷7=74, 6"
Replace-Text "698÷3=232, 2" "116÷3=38, 2"
Replace-Text "514÷9=57, 1" "221÷6=36, 5"

Replace-Text "999÷4=249, 3" "675÷8=84, 3"
Replace-Text "463÷3=154, 1" "137÷5=27, 2"
Replace-Text "347÷3=115, 2" "820÷8=102, 4"
Replace-Text "309÷2=154, 1" "248÷8=31, 0"
Replace-Text "587÷6=97, 5" "570÷8=71, 2"

Replace-Text "417÷3=139, 0" "369÷5=73, 4"
Replace-Text "355÷7=50, 5" "370÷7=52, 6"
Replace-Text "831÷3=277, 0" "959÷9=106, 5"
Replace-Text "670÷5=134, 0" "917÷4=229, 1"
Replace-Text "938÷4=234, 2" "303÷5=60, 3"

Replace-Text "219÷9=24, 3" "216÷5=43, 1"
Replace-Text "129÷3=43, 0" "386÷7=55, 1"
Replace-Text "469÷5=93, 4" "586÷9=65, 1"
Replace-Text "893÷4=223, 1" "292÷3=97, 1"
Replace-Text "930÷3=310, 0" "350÷4=87, 2"
